$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on the Price/Volume columns we touch so that
# Excel does not silently coerce strings like "304.50" or "0.999" into
# numbers (which would drop significant trailing zeros) or parse multi-dot
# values. This mirrors the source data, which stores these as inline strings.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.997.92'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.39%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.550.81'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.26%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '304.50'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.80%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '98.34'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +3.98%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.547'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.95%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.71'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +2.06%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0813'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.12%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.39%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.115'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +5.92%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.943.46'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.24%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.528.33'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.04%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.880'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.54%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.82'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +4.28%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.128.09'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.22%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.58'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +4.03%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0986'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.01%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.61'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.54%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '71.86'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.63%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '256.04'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.97%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.66%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.44%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '27.93'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -6.10%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.14%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.16'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.96%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '37.98'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +2.24%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.10'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.34%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.51%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '159.24'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.91%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.15%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.15'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.28%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0805'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.78%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.31'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -2.35%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '18.85'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +12.22%  '
$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.115'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.96%  '
$ws.Range("B39").Value = 'EnergySwap'
$ws.Range("C39").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '25.65'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +10.12%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.77%  '
$ws.Range("B41").Value = 'ApeXProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.09'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +31.51%  '
$ws.Range("B42").Value = 'NEARProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.43'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.89%  '
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.88'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.83%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.098.37'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.56%  '
$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0306'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -2.48%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.999'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.01%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '86.58'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.76%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.01'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.47%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '75.46'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +8.52%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.801.12'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.17%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '103.72'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.74%  '
